$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab: SCD0280 -> SCD0018 ---
$ws.Name = "SCD0018"

# --- Update the TC_ID (column B) for the three test-case rows ---
# The old Jira-style ID "DGS-295" is replaced everywhere by the new
# TC_ID scheme "SCD0018-003".
$ws.Range("B2").Value = "SCD0018-003"
$ws.Range("B3").Value = "SCD0018-003"
$ws.Range("B4").Value = "SCD0018-003"

# --- Column B needs to widen slightly to fit the new, longer TC_ID text ---
$ws.Range("B:B").ColumnWidth = 11.62

# --- Row 2 re-wraps with the new column width, so its height shrinks ---
$ws.Rows("2:2").RowHeight = 89.25

# --- Leave the final selection on C5, matching the saved workbook state ---
$ws.Range("C5").Select()
